# Auto-generated Excel COM-interop script
# Applies cell value updates to the cryptos worksheet per the target diff.
# Numeric-looking Price (column D) values are written with a leading apostrophe
# so Excel stores them as text (matching the source data's text-typed price column)
# instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.543.15'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '2.360.85'
$ws.Range('E3').Value = '  -3.88%  '
$ws.Range('D5').Value = '''541.78'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').Value = '''140.28'
$ws.Range('E6').Value = '  -3.63%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.533'
$ws.Range('E8').Value = '  -10.65%  '
$ws.Range('D9').Value = '2.359.38'
$ws.Range('E9').Value = '  -3.90%  '
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('D13').Value = '''0.341'
$ws.Range('E13').Value = '  -2.77%  '
$ws.Range('D14').Value = '''25.22'
$ws.Range('E14').Value = '  -3.01%  '
$ws.Range('D15').Value = '2.784.78'
$ws.Range('E15').Value = '  -3.91%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '60.556.38'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '''0.0000161'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '2.361.09'
$ws.Range('E18').Value = '  -3.79%  '
$ws.Range('D19').Value = '''10.52'
$ws.Range('E19').Value = '  -4.67%  '
$ws.Range('D20').Value = '''4.07'
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('D21').Value = '''313.42'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('E22').Value = '  -3.86%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('D25').Value = '''62.56'
$ws.Range('E25').Value = '  -1.35%  '
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').Value = '2.477.11'
$ws.Range('E27').Value = '  -4.04%  '
$ws.Range('D28').Value = '0.0₃0918'
$ws.Range('E28').Value = '  -5.88%  '
$ws.Range('D29').Value = '''7.64'
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').Value = '''513.24'
$ws.Range('E30').Value = '  -5.24%  '
$ws.Range('E31').Value = '  -4.77%  '
$ws.Range('D32').Value = '''7.90'
$ws.Range('E32').Value = '  -4.81%  '
$ws.Range('E33').Value = '  -4.21%  '
$ws.Range('E34').Value = '  -3.65%  '
$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '''5.41'
$ws.Range('E37').Value = '  -7.63%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '''4.61'
$ws.Range('E38').Value = '  -4.63%  '
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('D40').Value = '''17.89'
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('D43').Value = '''136.22'
$ws.Range('E43').Value = '  -6.64%  '
$ws.Range('D44').Value = '''40.19'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('D46').Value = '''137.62'
$ws.Range('E46').Value = '  -6.09%  '
$ws.Range('D47').Value = '''3.51'
$ws.Range('E47').Value = '  -1.09%  '
$ws.Range('D48').Value = '''20.11'
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('D49').Value = '''0.0513'
$ws.Range('E49').Value = '  -3.06%  '
$ws.Range('D50').Value = '''0.571'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').Value = '''0.0904'
$ws.Range('E51').Value = '  -3.65%  '

Write-Output "Applied updates to cryptos sheet"
